$wb = $excel.ActiveWorkbook

# Add the new worksheet after the current last sheet (ClusterExample),
# then rename it to ClusterTest.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ClusterTest"

# Header row
$ws.Cells.Item(1, 1).Value = "Name"
$ws.Cells.Item(1, 2).Value = "Age"
$ws.Cells.Item(1, 3).Value = "Income"

$data = @(
    @("Bob", 27, 70000),
    @("Michael", 29, 90000),
    @("Mohan", 29, 61000),
    @("Ismail", 28, 60000),
    @("Kory", 42, 150000),
    @("Gautam", 39, 155000),
    @("David", 41, 160000),
    @("Andrea", 38, 162000),
    @("Brad", 36, 156000),
    @("Angelina", 35, 130000),
    @("Donald", 37, 137000),
    @("Tom", 26, 45000),
    @("Arnold", 27, 48000),
    @("Jared", 28, 51000),
    @("Stark", 29, 49500),
    @("Ranbir", 32, 53000),
    @("Dipika", 40, 65000),
    @("Priyanka", 41, 63000),
    @("Nick", 43, 64000),
    @("Alisa", 39, 80000),
    @("Sid", 41, 82000),
    @("Tim", 39, 58000)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

# Match the author's final selection on the new sheet.
[void]$ws.Range("C25").Select()
